# Atualizacao de bases das ligas, do dia: 2024-01-29 as 18-52
#
# The underlying data for these match rows got refreshed/re-scraped; in this
# particular pass, 6 pairs of adjacent rows ended up with their betting-odds
# payload (everything from column B "id" through column AC "PL_AhUnder",
# except the row-sequence column A and the constant league columns C/D)
# swapped between the two rows of each pair. Column A (sequence number),
# C ("Div Original Name"), D ("Div"), and E (Date) are unaffected.
#
# Affected row pairs: (367,368) (382,383) (387,388) (443,444) (545,546) (595,596)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 367 ---
$ws.Cells.Item(367, 2).Value = 5172778   # B367
$ws.Cells.Item(367, 6).Value = "Lugano"   # F367
$ws.Cells.Item(367, 7).Value = "FC Zurich"   # G367
$ws.Cells.Item(367, 8).Value = 2   # H367
$ws.Cells.Item(367, 10).Value = "H"   # J367
$ws.Cells.Item(367, 11).Value = 2.1   # K367
$ws.Cells.Item(367, 12).Value = 3.6   # L367
$ws.Cells.Item(367, 13).Value = 3.3   # M367
$ws.Cells.Item(367, 14).Value = 2.7   # N367
$ws.Cells.Item(367, 15).Value = 3.3   # O367
$ws.Cells.Item(367, 16).Value = 2.625   # P367
$ws.Cells.Item(367, 17).Value = 0   # Q367
$ws.Cells.Item(367, 18).Value = 1.975   # R367
$ws.Cells.Item(367, 19).Value = 1.875   # S367
$ws.Cells.Item(367, 20).Value = 2.5   # T367
$ws.Cells.Item(367, 21).Value = 1.85   # U367
$ws.Cells.Item(367, 22).Value = 2   # V367
$ws.Cells.Item(367, 23).Value = 1.7   # W367
$ws.Cells.Item(367, 24).Value = -1   # X367
$ws.Cells.Item(367, 26).Value = 0.9750000000000001   # Z367
$ws.Cells.Item(367, 27).Value = -1   # AA367
$ws.Cells.Item(367, 29).Value = 1   # AC367

# --- Row 368 ---
$ws.Cells.Item(368, 2).Value = 5172777   # B368
$ws.Cells.Item(368, 6).Value = "Basel"   # F368
$ws.Cells.Item(368, 7).Value = "FC Sion"   # G368
$ws.Cells.Item(368, 8).Value = 0   # H368
$ws.Cells.Item(368, 10).Value = "D"   # J368
$ws.Cells.Item(368, 11).Value = 1.727   # K368
$ws.Cells.Item(368, 12).Value = 3.8   # L368
$ws.Cells.Item(368, 13).Value = 4.333   # M368
$ws.Cells.Item(368, 14).Value = 1.85   # N368
$ws.Cells.Item(368, 15).Value = 3.8   # O368
$ws.Cells.Item(368, 16).Value = 3.8   # P368
$ws.Cells.Item(368, 17).Value = -0.5   # Q368
$ws.Cells.Item(368, 18).Value = 1.85   # R368
$ws.Cells.Item(368, 19).Value = 2   # S368
$ws.Cells.Item(368, 20).Value = 3   # T368
$ws.Cells.Item(368, 21).Value = 2.05   # U368
$ws.Cells.Item(368, 22).Value = 1.8   # V368
$ws.Cells.Item(368, 23).Value = -1   # W368
$ws.Cells.Item(368, 24).Value = 2.8   # X368
$ws.Cells.Item(368, 26).Value = -1   # Z368
$ws.Cells.Item(368, 27).Value = 1   # AA368
$ws.Cells.Item(368, 29).Value = 0.8   # AC368

# --- Row 382 ---
$ws.Cells.Item(382, 2).Value = 5172789   # B382
$ws.Cells.Item(382, 6).Value = "FC Zurich"   # F382
$ws.Cells.Item(382, 7).Value = "St Gallen"   # G382
$ws.Cells.Item(382, 8).Value = 1   # H382
$ws.Cells.Item(382, 9).Value = 0   # I382
$ws.Cells.Item(382, 10).Value = "H"   # J382
$ws.Cells.Item(382, 11).Value = 2.625   # K382
$ws.Cells.Item(382, 12).Value = 3.4   # L382
$ws.Cells.Item(382, 13).Value = 2.5   # M382
$ws.Cells.Item(382, 14).Value = 2.25   # N382
$ws.Cells.Item(382, 15).Value = 3.6   # O382
$ws.Cells.Item(382, 16).Value = 3   # P382
$ws.Cells.Item(382, 17).Value = -0.25   # Q382
$ws.Cells.Item(382, 18).Value = 2   # R382
$ws.Cells.Item(382, 19).Value = 1.85   # S382
$ws.Cells.Item(382, 20).Value = 3   # T382
$ws.Cells.Item(382, 23).Value = 1.25   # W382
$ws.Cells.Item(382, 24).Value = -1   # X382
$ws.Cells.Item(382, 26).Value = 1   # Z382
$ws.Cells.Item(382, 27).Value = -1   # AA382
$ws.Cells.Item(382, 28).Value = -1   # AB382
$ws.Cells.Item(382, 29).Value = 0.8500000000000001   # AC382

# --- Row 383 ---
$ws.Cells.Item(383, 2).Value = 5172788   # B383
$ws.Cells.Item(383, 6).Value = "Servette"   # F383
$ws.Cells.Item(383, 7).Value = "FC Sion"   # G383
$ws.Cells.Item(383, 8).Value = 2   # H383
$ws.Cells.Item(383, 9).Value = 2   # I383
$ws.Cells.Item(383, 10).Value = "D"   # J383
$ws.Cells.Item(383, 11).Value = 2.1   # K383
$ws.Cells.Item(383, 12).Value = 3.3   # L383
$ws.Cells.Item(383, 13).Value = 3.4   # M383
$ws.Cells.Item(383, 14).Value = 2   # N383
$ws.Cells.Item(383, 15).Value = 3.5   # O383
$ws.Cells.Item(383, 16).Value = 3.75   # P383
$ws.Cells.Item(383, 17).Value = -0.5   # Q383
$ws.Cells.Item(383, 18).Value = 2.05   # R383
$ws.Cells.Item(383, 19).Value = 1.8   # S383
$ws.Cells.Item(383, 20).Value = 2.5   # T383
$ws.Cells.Item(383, 23).Value = -1   # W383
$ws.Cells.Item(383, 24).Value = 2.5   # X383
$ws.Cells.Item(383, 26).Value = -1   # Z383
$ws.Cells.Item(383, 27).Value = 0.8   # AA383
$ws.Cells.Item(383, 28).Value = 1   # AB383
$ws.Cells.Item(383, 29).Value = -1   # AC383

# --- Row 387 ---
$ws.Cells.Item(387, 2).Value = 5171748   # B387
$ws.Cells.Item(387, 6).Value = "Lucerne"   # F387
$ws.Cells.Item(387, 7).Value = "Young Boys"   # G387
$ws.Cells.Item(387, 9).Value = 1   # I387
$ws.Cells.Item(387, 10).Value = "D"   # J387
$ws.Cells.Item(387, 14).Value = 4.333   # N387
$ws.Cells.Item(387, 15).Value = 3.8   # O387
$ws.Cells.Item(387, 16).Value = 1.8   # P387
$ws.Cells.Item(387, 17).Value = 0.75   # Q387
$ws.Cells.Item(387, 18).Value = 1.825   # R387
$ws.Cells.Item(387, 19).Value = 2.025   # S387
$ws.Cells.Item(387, 20).Value = 3   # T387
$ws.Cells.Item(387, 21).Value = 2.025   # U387
$ws.Cells.Item(387, 22).Value = 1.825   # V387
$ws.Cells.Item(387, 23).Value = -1   # W387
$ws.Cells.Item(387, 24).Value = 2.8   # X387
$ws.Cells.Item(387, 26).Value = 0.825   # Z387
$ws.Cells.Item(387, 29).Value = 0.825   # AC387

# --- Row 388 ---
$ws.Cells.Item(388, 2).Value = 5173744   # B388
$ws.Cells.Item(388, 6).Value = "Winterthur"   # F388
$ws.Cells.Item(388, 7).Value = "Lugano"   # G388
$ws.Cells.Item(388, 9).Value = 0   # I388
$ws.Cells.Item(388, 10).Value = "H"   # J388
$ws.Cells.Item(388, 14).Value = 3   # N388
$ws.Cells.Item(388, 15).Value = 3.6   # O388
$ws.Cells.Item(388, 16).Value = 2.25   # P388
$ws.Cells.Item(388, 17).Value = 0.25   # Q388
$ws.Cells.Item(388, 18).Value = 1.9   # R388
$ws.Cells.Item(388, 19).Value = 1.95   # S388
$ws.Cells.Item(388, 20).Value = 2.5   # T388
$ws.Cells.Item(388, 21).Value = 1.875   # U388
$ws.Cells.Item(388, 22).Value = 1.975   # V388
$ws.Cells.Item(388, 23).Value = 2   # W388
$ws.Cells.Item(388, 24).Value = -1   # X388
$ws.Cells.Item(388, 26).Value = 0.8999999999999999   # Z388
$ws.Cells.Item(388, 29).Value = 0.9750000000000001   # AC388

# --- Row 443 ---
$ws.Cells.Item(443, 2).Value = 6401769   # B443
$ws.Cells.Item(443, 6).Value = "Lucerne"   # F443
$ws.Cells.Item(443, 7).Value = "Winterthur"   # G443
$ws.Cells.Item(443, 8).Value = 3   # H443
$ws.Cells.Item(443, 9).Value = 1   # I443
$ws.Cells.Item(443, 10).Value = "H"   # J443
$ws.Cells.Item(443, 11).Value = 1.5   # K443
$ws.Cells.Item(443, 12).Value = 4.333   # L443
$ws.Cells.Item(443, 13).Value = 5.5   # M443
$ws.Cells.Item(443, 14).Value = 1.444   # N443
$ws.Cells.Item(443, 15).Value = 4.75   # O443
$ws.Cells.Item(443, 16).Value = 6.5   # P443
$ws.Cells.Item(443, 17).Value = -1.25   # Q443
$ws.Cells.Item(443, 20).Value = 3   # T443
$ws.Cells.Item(443, 21).Value = 1.925   # U443
$ws.Cells.Item(443, 22).Value = 1.925   # V443
$ws.Cells.Item(443, 23).Value = 0.444   # W443
$ws.Cells.Item(443, 25).Value = -1   # Y443
$ws.Cells.Item(443, 26).Value = 0.9750000000000001   # Z443
$ws.Cells.Item(443, 27).Value = -1   # AA443
$ws.Cells.Item(443, 28).Value = 0.925   # AB443
$ws.Cells.Item(443, 29).Value = -1   # AC443

# --- Row 444 ---
$ws.Cells.Item(444, 2).Value = 6401770   # B444
$ws.Cells.Item(444, 6).Value = "FC Sion"   # F444
$ws.Cells.Item(444, 7).Value = "Basel"   # G444
$ws.Cells.Item(444, 8).Value = 1   # H444
$ws.Cells.Item(444, 9).Value = 2   # I444
$ws.Cells.Item(444, 10).Value = "A"   # J444
$ws.Cells.Item(444, 11).Value = 3.6   # K444
$ws.Cells.Item(444, 12).Value = 3.6   # L444
$ws.Cells.Item(444, 13).Value = 1.909   # M444
$ws.Cells.Item(444, 14).Value = 3.8   # N444
$ws.Cells.Item(444, 15).Value = 3.75   # O444
$ws.Cells.Item(444, 16).Value = 1.909   # P444
$ws.Cells.Item(444, 17).Value = 0.5   # Q444
$ws.Cells.Item(444, 20).Value = 2.75   # T444
$ws.Cells.Item(444, 21).Value = 1.95   # U444
$ws.Cells.Item(444, 22).Value = 1.9   # V444
$ws.Cells.Item(444, 23).Value = -1   # W444
$ws.Cells.Item(444, 25).Value = 0.909   # Y444
$ws.Cells.Item(444, 26).Value = -1   # Z444
$ws.Cells.Item(444, 27).Value = 0.875   # AA444
$ws.Cells.Item(444, 28).Value = 0.475   # AB444
$ws.Cells.Item(444, 29).Value = -0.5   # AC444

# --- Row 545 ---
$ws.Cells.Item(545, 2).Value = 6811248   # B545
$ws.Cells.Item(545, 6).Value = "Winterthur"   # F545
$ws.Cells.Item(545, 7).Value = "Young Boys"   # G545
$ws.Cells.Item(545, 8).Value = 1   # H545
$ws.Cells.Item(545, 9).Value = 4   # I545
$ws.Cells.Item(545, 10).Value = "A"   # J545
$ws.Cells.Item(545, 11).Value = 3.8   # K545
$ws.Cells.Item(545, 12).Value = 4   # L545
$ws.Cells.Item(545, 13).Value = 1.727   # M545
$ws.Cells.Item(545, 14).Value = 3.5   # N545
$ws.Cells.Item(545, 15).Value = 3.75   # O545
$ws.Cells.Item(545, 16).Value = 2   # P545
$ws.Cells.Item(545, 17).Value = 0.5   # Q545
$ws.Cells.Item(545, 18).Value = 1.825   # R545
$ws.Cells.Item(545, 19).Value = 2.025   # S545
$ws.Cells.Item(545, 20).Value = 3   # T545
$ws.Cells.Item(545, 21).Value = 2.025   # U545
$ws.Cells.Item(545, 22).Value = 1.825   # V545
$ws.Cells.Item(545, 23).Value = -1   # W545
$ws.Cells.Item(545, 25).Value = 1   # Y545
$ws.Cells.Item(545, 26).Value = -1   # Z545
$ws.Cells.Item(545, 27).Value = 1.025   # AA545
$ws.Cells.Item(545, 28).Value = 1.025   # AB545

# --- Row 546 ---
$ws.Cells.Item(546, 2).Value = 6811245   # B546
$ws.Cells.Item(546, 6).Value = "Lausanne Sports"   # F546
$ws.Cells.Item(546, 7).Value = "Lugano"   # G546
$ws.Cells.Item(546, 8).Value = 3   # H546
$ws.Cells.Item(546, 9).Value = 1   # I546
$ws.Cells.Item(546, 10).Value = "H"   # J546
$ws.Cells.Item(546, 11).Value = 2.25   # K546
$ws.Cells.Item(546, 12).Value = 3.4   # L546
$ws.Cells.Item(546, 13).Value = 2.875   # M546
$ws.Cells.Item(546, 14).Value = 2.375   # N546
$ws.Cells.Item(546, 15).Value = 3.4   # O546
$ws.Cells.Item(546, 16).Value = 3   # P546
$ws.Cells.Item(546, 17).Value = -0.25   # Q546
$ws.Cells.Item(546, 18).Value = 2.05   # R546
$ws.Cells.Item(546, 19).Value = 1.8   # S546
$ws.Cells.Item(546, 20).Value = 2.75   # T546
$ws.Cells.Item(546, 21).Value = 1.85   # U546
$ws.Cells.Item(546, 22).Value = 2   # V546
$ws.Cells.Item(546, 23).Value = 1.375   # W546
$ws.Cells.Item(546, 25).Value = -1   # Y546
$ws.Cells.Item(546, 26).Value = 1.05   # Z546
$ws.Cells.Item(546, 27).Value = -1   # AA546
$ws.Cells.Item(546, 28).Value = 0.8500000000000001   # AB546

# --- Row 595 ---
$ws.Cells.Item(595, 2).Value = 6810780   # B595
$ws.Cells.Item(595, 6).Value = "Young Boys"   # F595
$ws.Cells.Item(595, 7).Value = "Yverdon Sport FC"   # G595
$ws.Cells.Item(595, 11).Value = 1.363   # K595
$ws.Cells.Item(595, 12).Value = 5   # L595
$ws.Cells.Item(595, 13).Value = 8   # M595
$ws.Cells.Item(595, 14).Value = 1.4   # N595
$ws.Cells.Item(595, 15).Value = 4.75   # O595
$ws.Cells.Item(595, 16).Value = 7.5   # P595
$ws.Cells.Item(595, 17).Value = -1.5   # Q595
$ws.Cells.Item(595, 18).Value = 2.025   # R595
$ws.Cells.Item(595, 19).Value = 1.825   # S595
$ws.Cells.Item(595, 20).Value = 3   # T595
$ws.Cells.Item(595, 21).Value = 1.875   # U595
$ws.Cells.Item(595, 22).Value = 1.975   # V595

# --- Row 596 ---
$ws.Cells.Item(596, 2).Value = 6811278   # B596
$ws.Cells.Item(596, 6).Value = "Winterthur"   # F596
$ws.Cells.Item(596, 7).Value = "Basel"   # G596
$ws.Cells.Item(596, 11).Value = 2.3   # K596
$ws.Cells.Item(596, 12).Value = 3.6   # L596
$ws.Cells.Item(596, 13).Value = 2.875   # M596
$ws.Cells.Item(596, 14).Value = 2.4   # N596
$ws.Cells.Item(596, 15).Value = 3.6   # O596
$ws.Cells.Item(596, 16).Value = 2.8   # P596
$ws.Cells.Item(596, 17).Value = -0.25   # Q596
$ws.Cells.Item(596, 18).Value = 2.05   # R596
$ws.Cells.Item(596, 19).Value = 1.8   # S596
$ws.Cells.Item(596, 20).Value = 2.75   # T596
$ws.Cells.Item(596, 21).Value = 2   # U596
$ws.Cells.Item(596, 22).Value = 1.85   # V596

